$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add row 5 to Sheet1
$ws1.Range("A5").Value = "raju"
$ws1.Range("B5").Value = "SSS"
$ws1.Range("C5").Value = "vashi"

# Update selection on Sheet1
$null = $ws1.Range("E21").Select()

# Add Sheet2 after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "aazz"
$ws2.Range("B1").Value = "terna"
$ws2.Range("C1").Value = "nerul"
$ws2.Range("D1").Value = 1

$ws2.Range("A2").Value = "bb"
$ws2.Range("B2").Value = "MGM"
$ws2.Range("C2").Value = "vashi"
$ws2.Range("D2").Value = 2

$ws2.Range("A3").Value = "cc"
$ws2.Range("B3").Value = "SS"
$ws2.Range("C3").Value = "seawoods"
$ws2.Range("D3").Value = 3

$ws2.Range("A4").Value = "dd"
$ws2.Range("B4").Value = "agnel"
$ws2.Range("C4").Value = "parsik"
$ws2.Range("D4").Value = 4

$ws2.Range("A5").Value = "ee"
$ws2.Range("B5").Value = "SSS"
$ws2.Range("C5").Value = "vashi"
$ws2.Range("D5").Value = 5

$null = $ws2.Range("D5").Select()

$null = $ws1.Activate()
